function Set-TextCell {
    param($ws, $row, $col, $val)
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-TextCell $ws 2 4 "41.347.06"
Set-TextCell $ws 2 5 "  +1.54%  "

# Row 3
Set-TextCell $ws 3 4 "2.184.00"
Set-TextCell $ws 3 5 "  -0.07%  "

# Row 4
Set-TextCell $ws 4 5 "  -0.12%  "

# Row 5
Set-TextCell $ws 5 4 "254.93"
Set-TextCell $ws 5 5 "  +5.78%  "

# Row 6
Set-TextCell $ws 6 4 "0.627"
Set-TextCell $ws 6 5 "  +1.14%  "

# Row 7
Set-TextCell $ws 7 4 "68.17"
Set-TextCell $ws 7 5 "  -0.24%  "

# Row 8
Set-TextCell $ws 8 5 "  -0.14%  "

# Row 9
Set-TextCell $ws 9 4 "0.576"
Set-TextCell $ws 9 5 "  +7.50%  "

# Row 10
Set-TextCell $ws 10 4 "37.79"
Set-TextCell $ws 10 5 "  +4.76%  "

# Row 11
Set-TextCell $ws 11 4 "59.05"
Set-TextCell $ws 11 5 "  +2.70%  "

# Row 12
Set-TextCell $ws 12 4 "0.0935"
Set-TextCell $ws 12 5 "  -0.47%  "

# Row 13
Set-TextCell $ws 13 4 "7.14"
Set-TextCell $ws 13 5 "  +8.99%  "

# Row 14
Set-TextCell $ws 14 5 "  +0.85%  "

# Row 15
Set-TextCell $ws 15 4 "2.499.55"
Set-TextCell $ws 15 5 "  -0.49%  "

# Row 16
Set-TextCell $ws 16 4 "0.872"
Set-TextCell $ws 16 5 "  +5.43%  "

# Row 17
Set-TextCell $ws 17 4 "14.46"
Set-TextCell $ws 17 5 "  -0.79%  "

# Row 18
Set-TextCell $ws 18 4 "2.180.85"
Set-TextCell $ws 18 5 "  -0.31%  "

# Row 19
Set-TextCell $ws 19 4 "41.199.29"
Set-TextCell $ws 19 5 "  +1.38%  "

# Row 20
Set-TextCell $ws 20 4 "0.0₃0953"
Set-TextCell $ws 20 5 "  +1.84%  "

# Row 21
Set-TextCell $ws 21 4 "6.17"
Set-TextCell $ws 21 5 "  +2.44%  "

# Row 22
Set-TextCell $ws 22 4 "71.81"
Set-TextCell $ws 22 5 "  -0.58%  "

# Row 23
Set-TextCell $ws 23 4 "231.85"
Set-TextCell $ws 23 5 "  +1.56%  "

# Row 24
Set-TextCell $ws 24 4 "2.03"

# Row 25
Set-TextCell $ws 25 4 "3.95"
Set-TextCell $ws 25 5 "  +10.31%  "

# Row 26
Set-TextCell $ws 26 4 "11.74"
Set-TextCell $ws 26 5 "  +22.16%  "

# Row 27
Set-TextCell $ws 27 5 "  -0.07%  "

# Row 28
Set-TextCell $ws 28 4 "2.54"
Set-TextCell $ws 28 5 "  +6.51%  "

# Row 29
Set-TextCell $ws 29 5 "  +0.29%  "

# Row 30
Set-TextCell $ws 30 4 "168.40"
Set-TextCell $ws 30 5 "  +0.05%  "

# Row 31
Set-TextCell $ws 31 4 "20.64"
Set-TextCell $ws 31 5 "  +2.60%  "

# Row 32
Set-TextCell $ws 32 4 "0.117"
Set-TextCell $ws 32 5 "  -0.20%  "

# Row 33
Set-TextCell $ws 33 4 "0.0749"
Set-TextCell $ws 33 5 "  +7.73%  "

# Row 34
Set-TextCell $ws 34 5 "  +0.65%  "

# Row 35
Set-TextCell $ws 35 4 "5.47"
Set-TextCell $ws 35 5 "  +7.24%  "

# Row 36
Set-TextCell $ws 36 4 "26.35"
Set-TextCell $ws 36 5 "  +12.39%  "

# Row 37
Set-TextCell $ws 37 2 "Filecoin"
Set-TextCell $ws 37 3 "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell $ws 37 4 "4.63"
Set-TextCell $ws 37 5 "  +1.69%  "

# Row 38
Set-TextCell $ws 38 2 "RenderToken"
Set-TextCell $ws 38 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws 38 4 "4.14"
Set-TextCell $ws 38 5 "  +8.75%  "

# Row 39
Set-TextCell $ws 39 4 "0.0301"
Set-TextCell $ws 39 5 "  +12.51%  "

# Row 40
Set-TextCell $ws 40 4 "2.20"
Set-TextCell $ws 40 5 "  -2.27%  "

# Row 41
Set-TextCell $ws 41 4 "12.50"
Set-TextCell $ws 41 5 "  +22.32%  "

# Row 42
Set-TextCell $ws 42 4 "5.67"
Set-TextCell $ws 42 5 "  -1.18%  "

# Row 43
Set-TextCell $ws 43 4 "64.22"
Set-TextCell $ws 43 5 "  +3.84%  "

# Row 44
Set-TextCell $ws 44 4 "5.06"
Set-TextCell $ws 44 5 "  +4.98%  "

# Row 45
Set-TextCell $ws 45 4 "0.201"
Set-TextCell $ws 45 5 "  +6.58%  "

# Row 46
Set-TextCell $ws 46 5 "  +0.95%  "

# Row 47
Set-TextCell $ws 47 5 "  +3.54%  "

# Row 48
Set-TextCell $ws 48 4 "1.00"
Set-TextCell $ws 48 5 "  +0.31%  "

# Row 49
Set-TextCell $ws 49 4 "1.14"
Set-TextCell $ws 49 5 "  +5.07%  "

# Row 50
Set-TextCell $ws 50 5 "  +1.57%  "

# Row 51
Set-TextCell $ws 51 4 "4.26"
Set-TextCell $ws 51 5 "  -3.69%  "
